$wb = $excel.ActiveWorkbook

# 1. Duplicate the "Sedan_Hamba" sheet (same style/template family as the new
#    FSAE_Achilles tab) and move the copy to the end of the workbook.
$src = $wb.Worksheets.Item("Sedan_Hamba")
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy([System.Reflection.Missing]::Value, $last)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "FSAE_Achilles"

# 2. Update the body data on the new sheet.
$newSheet.Range("H3").Value = "FSAE_Achilles"

$newSheet.Range("F6").Value = -1.53
$newSheet.Range("G6").Value = 0
$newSheet.Range("H6").Value = 0

$newSheet.Range("F7").Value = -0.8
$newSheet.Range("G7").Value = 0
$newSheet.Range("H7").Value = 0.289

$newSheet.Range("F8").Value = -1
$newSheet.Range("G8").Value = 0
$newSheet.Range("H8").Value = 0

$newSheet.Range("F9").Value = 0.25
$newSheet.Range("G9").Value = 0
$newSheet.Range("H9").Value = 0.403

$newSheet.Range("F10").Value = -1.75
$newSheet.Range("G10").Value = 0
$newSheet.Range("H10").Value = 0.403

$newSheet.Range("H11").Formula = "=0.619*2+0.2"

$newSheet.Range("H12").Value = 165

$newSheet.Range("F13").Value = 43
$newSheet.Range("G13").Value = 192
$newSheet.Range("H13").Value = 206

# 3. Fix up the original Sedan_Hamba sheet's remembered cursor position.
$src.Range("E22").Select()

# 4. Fix up the Sedan_HambaLG sheet's remembered cursor position (it was the
#    previously-active tab; the new sheet takes that role now).
$ws2 = $wb.Worksheets.Item("Sedan_HambaLG")
$ws2.Activate()
$ws2.Range("H12").Select()

# 5. Make the new sheet the active / selected tab, with its own cursor
#    position.
$newSheet.Activate()
$newSheet.Range("G27").Select()

Write-Output "done"
